# PSI_inputs.xlsx edit: "Changing residual embedment calculation method and
# some troubleshooting on the UD method weighting calcs"
#
# 1. On the "Inputs" sheet, the residual-embedment suction switch (B11) is
#    switched on from 0 to 2.
# 2. The "Automated Fit" distribution choice list ("N"/"LN") that lived in
#    column G alongside the "Base Variable" helper column in H is dropped;
#    the "Base Variable" values move left from H into G (rows 21:46).
# 3. The weighting values used for the UD method (row 37, "Lat_brk_weighting")
#    are updated.
# 4. The view is scrolled down and a new cell is selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inputs")

# --- 1. Residual embedment: turn on suction switch ---------------------
$ws.Range("B11").Value = 2

# --- 2. Drop the old "N"/"LN" dropdown column (G21:G46) and shift the ---
#        "Base Variable" helper column (H21:H46) one column left into G.
$ws.Range("G21:G46").ClearContents()
$ws.Range("H21:H46").Cut($ws.Range("G21:G46"))

# --- 3. Troubleshooting the UD method weighting calcs (row 37) ---------
$ws.Range("B37").Value = 1
$ws.Range("C37").Value = 1.5
$ws.Range("D37").Value = 2

# --- 4. Scroll / selection bookkeeping ----------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("P27").Select()
